$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Helper marker used to work around a boundary quirk when placing a
# zero-length bookmark exactly at the end of a paragraph's text (right
# before the paragraph mark): we temporarily append a marker string,
# anchor the bookmark just before it, then remove the marker again.
# -----------------------------------------------------------------------
$marker = "@@MARK@@"

# =========================================================================
# 1) First paragraph ("Predictive modeling"):
#    - insert a leading space run
#    - move the "_GoBack" bookmark here (it currently sits inside the
#      'Elastic net' paragraph further down)
# =========================================================================

$rng = $d.Range(0, 0)
$rng.Find.Execute("Predictive modeling") | Out-Null
$rng.InsertBefore(" ")

# remove the old _GoBack bookmark (currently in the Elastic net paragraph)
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()

# re-find "Predictive modeling" (position may have shifted) and append the
# temp marker right after it
$rng = $d.Range(0, 0)
$rng.Find.Execute("Predictive modeling") | Out-Null
$rng.InsertAfter($marker)

# anchor the bookmark collapsed right before the marker (i.e. right after
# "Predictive modeling")
$mrng = $d.Range(0, 0)
$mrng.Find.Execute($marker) | Out-Null
$bmRng = $d.Range($mrng.Start, $mrng.Start)
$d.Bookmarks.Add("_GoBack", $bmRng)

# remove the temp marker again
$mrng = $d.Range(0, 0)
$mrng.Find.Execute($marker) | Out-Null
$mrng.Delete()

# =========================================================================
# 2) Remove the empty paragraph that used to separate "Predictive modeling"
#    from the "Support Vector Machines..." paragraph.
# =========================================================================
$d.Paragraphs(2).Range.Delete()

# =========================================================================
# 3) Drop the trailing sentence from the SVM paragraph.
# =========================================================================
$d.Content.Find.Execute(
    "SVM is a statistical model based where the high dimensional feature space contains the input data that are mapped by a " + [char]0x2018 + "kernel function" + [char]0x2019 + " (nonlinear). ",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# =========================================================================
# 4) Shorten the "To facilitate computation efficiency..." lead-in and
#    change "determine which" -> "prefilter".
# =========================================================================
$d.Content.Find.Execute(
    "To facilitate computation efficiency, only CpGs with variability were chosen to train the models on. We chose to use an arbitrary threshold to determine which CpGs to retain.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "We chose to use an arbitrary threshold to prefilter CpGs to retain.", 2) | Out-Null

# =========================================================================
# 5) Re-word the Elastic net / SVM tuning-parameter sentence.
# =========================================================================
$d.Content.Find.Execute(
    "used for Elastic net, while 5 values of C for SVM were tested.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "tested for elastic net.", 2) | Out-Null

Write-Output "Done"
